$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 210.1399993896484
$ws.Range("E2").Value = 200.2799987792969
$ws.Range("F2").Value = 223.9499969482422
$ws.Range("G2").Value = 199.6000061035156
$ws.Range("H2").Value = 278228356
$ws.Range("I2").Value = "BIDU"

$ws.Range("D3").Value = 200.5899963378907
$ws.Range("E3").Value = 172.6600036621094
$ws.Range("F3").Value = 210
$ws.Range("G3").Value = 162
$ws.Range("H3").Value = 278228356
$ws.Range("I3").Value = "BIDU"

$ws.Range("D4").Value = 138.3200073242188
$ws.Range("E4").Value = 187.4700012207031
$ws.Range("F4").Value = 189.5500030517578
$ws.Range("G4").Value = 135.3099975585938
$ws.Range("H4").Value = 278228356
$ws.Range("I4").Value = "BIDU"

$ws.Range("D5").Value = 181.2899932861328
$ws.Range("E5").Value = 163.2700042724609
$ws.Range("F5").Value = 189.8399963378907
$ws.Range("G5").Value = 154.5299987792969
$ws.Range("H5").Value = 278228356
$ws.Range("I5").Value = "BIDU"

$ws.Range("D6").Value = 188.7899932861328
$ws.Range("E6").Value = 194.3000030517578
$ws.Range("F6").Value = 201
$ws.Range("G6").Value = 181.4400024414062
$ws.Range("H6").Value = 278228356
$ws.Range("I6").Value = "BIDU"

$ws.Range("D7").Value = 165.2400054931641
$ws.Range("E7").Value = 159.6000061035156
$ws.Range("F7").Value = 172.7599945068359
$ws.Range("G7").Value = 155.2799987792969
$ws.Range("H7").Value = 278228356
$ws.Range("I7").Value = "BIDU"

$ws.Range("D8").Value = 182.8399963378907
$ws.Range("E8").Value = 176.8600006103516
$ws.Range("F8").Value = 187.2400054931641
$ws.Range("G8").Value = 171.8899993896484
$ws.Range("H8").Value = 278228356
$ws.Range("I8").Value = "BIDU"

$ws.Range("D9").Value = 166.3000030517578
$ws.Range("E9").Value = 175.0700073242188
$ws.Range("F9").Value = 183
$ws.Range("G9").Value = 165.8200073242188
$ws.Range("H9").Value = 278228356
$ws.Range("I9").Value = "BIDU"

$ws.Range("D10").Value = 173
$ws.Range("E10").Value = 180.229995727539
$ws.Range("F10").Value = 188.6000061035156
$ws.Range("G10").Value = 171.1699981689453
$ws.Range("H10").Value = 278228356
$ws.Range("I10").Value = "BIDU"

$ws.Range("D11").Value = 179.8999938964844
$ws.Range("E11").Value = 226.3500061035156
$ws.Range("F11").Value = 230
$ws.Range("G11").Value = 179.6300048828125
$ws.Range("H11").Value = 278228356
$ws.Range("I11").Value = "BIDU"

$ws.Range("D12").Value = 248.6199951171875
$ws.Range("E12").Value = 243.9400024414062
$ws.Range("F12").Value = 274.9700012207031
$ws.Range("G12").Value = 232.4400024414062
$ws.Range("H12").Value = 278228356
$ws.Range("I12").Value = "BIDU"

$ws.Range("D13").Value = 236.4900054931641
$ws.Range("E13").Value = 246.9199981689453
$ws.Range("F13").Value = 265.7699890136719
$ws.Range("G13").Value = 235.259994506836
$ws.Range("H13").Value = 278228356
$ws.Range("I13").Value = "BIDU"

$ws.Range("D14").Value = 220.8399963378907
$ws.Range("E14").Value = 250.8999938964844
$ws.Range("F14").Value = 257.2999877929688
$ws.Range("G14").Value = 213.5599975585937
$ws.Range("H14").Value = 278228356
$ws.Range("I14").Value = "BIDU"

$ws.Range("D15").Value = 238.0800018310547
$ws.Range("E15").Value = 247.1799926757812
$ws.Range("F15").Value = 274
$ws.Range("G15").Value = 236.2799987792969
$ws.Range("H15").Value = 278228356
$ws.Range("I15").Value = "BIDU"

$ws.Range("D16").Value = 230.8099975585937
$ws.Range("E16").Value = 190.0599975585937
$ws.Range("F16").Value = 231.1699981689453
$ws.Range("G16").Value = 177.8000030517578
$ws.Range("H16").Value = 278228356
$ws.Range("I16").Value = "BIDU"

$ws.Range("D17").Value = 156.1799926757812
$ws.Range("E17").Value = 172.6300048828125
$ws.Range("F17").Value = 175.0899963378906
$ws.Range("G17").Value = 153.7799987792969
$ws.Range("H17").Value = 278228356
$ws.Range("I17").Value = "BIDU"

$ws.Range("D18").Value = 168.0700073242188
$ws.Range("E18").Value = 166.2299957275391
$ws.Range("F18").Value = 186.2200012207031
$ws.Range("G18").Value = 162.1999969482422
$ws.Range("H18").Value = 278228356
$ws.Range("I18").Value = "BIDU"

$ws.Range("D19").Value = 121.4899978637695
$ws.Range("E19").Value = 111.6999969482422
$ws.Range("F19").Value = 121.8000030517578
$ws.Range("G19").Value = 110.75
$ws.Range("H19").Value = 278228356
$ws.Range("I19").Value = "BIDU"

$ws.Range("D20").Value = 102.8000030517578
$ws.Range("E20").Value = 101.8499984741211
$ws.Range("F20").Value = 109.3300018310547
$ws.Range("G20").Value = 98.1999969482422
$ws.Range("H20").Value = 278228356
$ws.Range("I20").Value = "BIDU"

$ws.Range("D21").Value = 129.4900054931641
$ws.Range("E21").Value = 123.5599975585938
$ws.Range("F21").Value = 147.3800048828125
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = 278228356
$ws.Range("I21").Value = "BIDU"

$ws.Range("D22").Value = 98.68000030517578
$ws.Range("E22").Value = 100.9300003051758
$ws.Range("F22").Value = 107.6699981689453
$ws.Range("G22").Value = 94.90000152587891
$ws.Range("H22").Value = 278228356
$ws.Range("I22").Value = "BIDU"

$ws.Range("D23").Value = 119.8899993896484
$ws.Range("E23").Value = 119.4000015258789
$ws.Range("F23").Value = 135.4400024414062
$ws.Range("G23").Value = 114.75
$ws.Range("H23").Value = 278228356
$ws.Range("I23").Value = "BIDU"

$ws.Range("D24").Value = 126.5299987792969
$ws.Range("E24").Value = 133.0500030517578
$ws.Range("F24").Value = 138.9799957275391
$ws.Range("G24").Value = 123.5
$ws.Range("H24").Value = 278228356
$ws.Range("I24").Value = "BIDU"

$ws.Range("D25").Value = 219.509994506836
$ws.Range("E25").Value = 235.0200042724609
$ws.Range("F25").Value = 264.9400024414062
$ws.Range("G25").Value = 203.6199951171875
$ws.Range("H25").Value = 278228356
$ws.Range("I25").Value = "BIDU"

$ws.Range("D26").Value = 224
$ws.Range("E26").Value = 210.3300018310547
$ws.Range("F26").Value = 228.5299987792969
$ws.Range("G26").Value = 204.8800048828125
$ws.Range("H26").Value = 278228356
$ws.Range("I26").Value = "BIDU"

$ws.Range("D27").Value = 203.9100036621093
$ws.Range("E27").Value = 164.0099945068359
$ws.Range("F27").Value = 205.4100036621093
$ws.Range("G27").Value = 153.1399993896484
$ws.Range("H27").Value = 278228356
$ws.Range("I27").Value = "BIDU"

$ws.Range("D28").Value = 153.5899963378906
$ws.Range("E28").Value = 162.2400054931641
$ws.Range("F28").Value = 182.6000061035156
$ws.Range("G28").Value = 144.5099945068359
$ws.Range("H28").Value = 278228356
$ws.Range("I28").Value = "BIDU"

$ws.Range("D29").Value = 148.9100036621094
$ws.Range("E29").Value = 159.7400054931641
$ws.Range("F29").Value = 165.0200042724609
$ws.Range("G29").Value = 139.1000061035156
$ws.Range("H29").Value = 278228356
$ws.Range("I29").Value = "BIDU"

$ws.Range("D30").Value = 144.0099945068359
$ws.Range("E30").Value = 124.1699981689453
$ws.Range("F30").Value = 154.2899932861328
$ws.Range("G30").Value = 110.3000030517578
$ws.Range("H30").Value = 278228356
$ws.Range("I30").Value = "BIDU"

$ws.Range("D31").Value = 148.0899963378906
$ws.Range("E31").Value = 136.5700073242188
$ws.Range("F31").Value = 156.6900024414062
$ws.Range("G31").Value = 131.0500030517578
$ws.Range("H31").Value = 278228356
$ws.Range("I31").Value = "BIDU"

$ws.Range("D32").Value = 117.1999969482422
$ws.Range("E32").Value = 76.56999969482422
$ws.Range("F32").Value = 125.2300033569336
$ws.Range("G32").Value = 73.58000183105469
$ws.Range("H32").Value = 278228356
$ws.Range("I32").Value = "BIDU"

$ws.Range("D33").Value = 118.8600006103516
$ws.Range("E33").Value = 134.6799926757812
$ws.Range("F33").Value = 142.4100036621094
$ws.Range("G33").Value = 118.0400009155273
$ws.Range("H33").Value = 278228356
$ws.Range("I33").Value = "BIDU"

$ws.Range("D34").Value = 150
$ws.Range("E34").Value = 120.6100006103516
$ws.Range("F34").Value = 151.0899963378906
$ws.Range("G34").Value = 116.370002746582
$ws.Range("H34").Value = 278228356
$ws.Range("I34").Value = "BIDU"

$ws.Range("D35").Value = 142.1499938964844
$ws.Range("E35").Value = 155.9900054931641
$ws.Range("F35").Value = 156.9799957275391
$ws.Range("G35").Value = 137.4199981689453
$ws.Range("H35").Value = 278228356
$ws.Range("I35").Value = "BIDU"

$ws.Range("D36").Value = 133.3200073242188
$ws.Range("E36").Value = 105
$ws.Range("F36").Value = 135.8500061035156
$ws.Range("G36").Value = 103.370002746582
$ws.Range("H36").Value = 278228356
$ws.Range("I36").Value = "BIDU"

$ws.Range("D37").Value = 116.8899993896484
$ws.Range("E37").Value = 105.3099975585938
$ws.Range("F37").Value = 119.75
$ws.Range("G37").Value = 97.51000213623048
$ws.Range("H37").Value = 278228356
$ws.Range("I37").Value = "BIDU"

$ws.Range("D38").Value = 107
$ws.Range("E38").Value = 103.4000015258789
$ws.Range("F38").Value = 109.9100036621094
$ws.Range("G38").Value = 94.25
$ws.Range("H38").Value = 278228356
$ws.Range("I38").Value = "BIDU"

$ws.Range("D39").Value = 86.61000061035156
$ws.Range("E39").Value = 88.56999969482422
$ws.Range("F39").Value = 104.6999969482422
$ws.Range("G39").Value = 85.08000183105469
$ws.Range("H39").Value = 278228356
$ws.Range("I39").Value = "BIDU"

$ws.Range("D40").Value = 105.75
$ws.Range("E40").Value = 91.23000335693359
$ws.Range("F40").Value = 116.25
$ws.Range("G40").Value = 88.58000183105469
$ws.Range("H40").Value = 278228356
$ws.Range("I40").Value = "BIDU"

$ws.Range("D41").Value = 83.16999816894531
$ws.Range("E41").Value = 90.59999847412109
$ws.Range("F41").Value = 96.18000030517578
$ws.Range("G41").Value = 77.19000244140625
$ws.Range("H41").Value = 278228356
$ws.Range("I41").Value = "BIDU"

$ws.Range("D42").Value = 91.09999847412109
$ws.Range("E42").Value = 87.81999969482422
$ws.Range("F42").Value = 92.70999908447266
$ws.Range("G42").Value = 74.70999908447266
$ws.Range("H42").Value = 278228356
$ws.Range("I42").Value = "BIDU"

$ws.Range("D43").Value = 85.68000030517578
$ws.Range("E43").Value = 87.87000274658203
$ws.Range("F43").Value = 94.5
$ws.Range("G43").Value = 87.87000274658203
$ws.Range("I43").Value = "BIDU"
